$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "76.012.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.58%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.909.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.40%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "198.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "598.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.75%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.551"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.10%  "

$ws.Range("E9").Value = "  -0.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.907.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.42%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.425"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +14.08%  "

$ws.Range("E12").Value = "  -1.17%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.74%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.433.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "75.869.65"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.34%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000190"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.49%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.903.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.84%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.73%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "377.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("E22").Value = "  +0.96%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.79%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.054.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.02%  "

$ws.Range("E27").Value = "  -0.81%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.79%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000109"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "501.64"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.67%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.73%  "

$ws.Range("E34").Value = "  -1.11%  "

$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.12%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.64%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.70"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.15%  "

$ws.Range("E39").Value = "  -6.76%  "

$ws.Range("E40").Value = "  -0.08%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "180.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.342"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.04%  "

$ws.Range("E44").Value = "  -1.97%  "

$ws.Range("E45").Value = "  +6.30%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.22%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.78%  "

$ws.Range("E48").Value = "  -2.26%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.576"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.657"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.29%  "
